# Swap the presentation's colour theme from the "Integral" palette to the
# standard "Office Theme" palette (the 12 theme colours), matching the
# colour values that the Office Theme uses. This is applied through the
# presentation's ThemeColorScheme, which is backed by the deck's theme
# part (ppt/theme/theme2.xml).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (R, G, B) using the standard 12-slot MsoThemeColorSchemeIndex
# ordering: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    @(0, 0, 0),         # dk1
    @(255, 255, 255),   # lt1
    @(68, 84, 106),      # dk2
    @(231, 230, 230),    # lt2
    @(91, 155, 213),     # accent1
    @(237, 125, 49),     # accent2
    @(165, 165, 165),    # accent3
    @(255, 192, 0),      # accent4
    @(68, 114, 196),     # accent5
    @(112, 173, 71),     # accent6
    @(5, 99, 193),       # hlink
    @(149, 79, 114)      # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $r = $rgb[0]
    $g = $rgb[1]
    $b = $rgb[2]
    $tcs.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
